$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": insert a new data row 47 (pushing the footnote row to 48)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(47).Insert()
$wsAll.Range("A47").Value = 43975
$wsAll.Range("B47").Value = 285
$wsAll.Range("C47").Value = 282
$wsAll.Range("D47").Value = 32
$wsAll.Range("E47").Value = 28
$wsAll.Range("F47").Value = 4
$wsAll.Range("G47").Value = 12
$wsAll.Range("H47").Value = 238
$wsAll.Range("H47").Select()

# ---------------------------------------------------------------------------
# Sheet "kobe": insert a new data row 102 (pushing the footnote row to 103)
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows.Item(102).Insert()
$wsKobe.Range("A102").Value = 43975
$wsKobe.Range("B102").Value = 0
$wsKobe.Range("C102").Value = 3021
$wsKobe.Range("D102").Value = 0
$wsKobe.Range("E102").Value = 285
$wsKobe.Range("F102").Value = 27
$wsKobe.Range("G102").Value = 24
$wsKobe.Range("H102").Value = 3
$wsKobe.Range("I102").Value = 12
$wsKobe.Range("J102").Value = 229
$wsKobe.Range("G103").Select()

# ---------------------------------------------------------------------------
# Sheet "other": insert a new data row 77 (pushing the footnote row to 78)
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(77).Insert()
$wsOther.Range("A77").Value = 43975
$wsOther.Range("B77").Value = 0
$wsOther.Range("C77").Value = 14
$wsOther.Range("D77").Value = 5
$wsOther.Range("E77").Value = 4
$wsOther.Range("F77").Value = 1
$wsOther.Range("G77").Value = 0
$wsOther.Range("H77").Value = 9
$wsOther.Range("H77").Select()

# Restore the originally active sheet/tab so tabSelected stays on "all".
$wsAll.Activate()
